$d = $word.ActiveDocument

# Make sure edits are applied cleanly (not as tracked-change insert/delete pairs).
$d.TrackRevisions = $false

# 1) Update the timestamp text near the top of the document.
$d.Content.Find.Execute("11/18/2023 4:13:15 PM", $true, $false, $false, $false, $false, $true, 1, $false, "11/18/2023 4:55:21 PM", 2) | Out-Null

# 2) Insert a new "PREVENTION SECURITY SYSTEM: ANY GUN MUZZLE;" paragraph
#    right after the "...ALL DIGITAL WEAPON TYPES;" paragraph (i.e. right
#    before the existing "...ANY LIGHTNING STICK;" paragraph).
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("LIGHTNING", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'ANY LIGHTNING STICK' paragraph to anchor the new paragraph."
}
$lightningPara = $findRange.Paragraphs(1)
$lightningPara.Range.InsertParagraphBefore() | Out-Null

# Re-resolve the freshly inserted (still empty) paragraph and fill it with
# the same OOXML structure used by its sibling "PREVENTION SECURITY SYSTEM"
# entries, but describing "ANY GUN MUZZLE".
$findRange2 = $d.Content.Duplicate
$findRange2.Find.Execute("LIGHTNING", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lightningPara2 = $findRange2.Paragraphs(1)
$newPara = $lightningPara2.Previous()
$newRange = $newPara.Range

$newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:ind w:left="720"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>PREVENTION SECURITY SYSTEM</w:t></w:r>
<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">ANY </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">GUN </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>MUZZLE</w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$newRange.InsertXML($newParaXml) | Out-Null
